# "update chức năng thêm lớp mới" — reshape the "1A1" header row:
#  - drop "Mã số học sinh" and "Trạng thái"
#  - replace "Fist Name"/"Last Name" with "Họ"/"Tên"
#  - move "Năm học"/"Khối"/"Lớp" from the tail (AA:AC) to the front (B:D)
#  - everything else keeps its relative order, just shifted into the new layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row, in final column order (A1 .. AA1)
$headers = @(
    "STT",
    "Năm học",
    "Khối",
    "Lớp",
    "Họ",
    "Tên",
    "Năm sinh",
    "Giới tính",
    "Dân tộc",
    "Ngày vào trường",
    "Số điện thoại",
    "Địa chỉ",
    "Cha",
    "Mẹ",
    "Quan hệ khác",
    "Họ tên cha",
    "Năm sinh cha",
    "Số điện thoại cha",
    "Nghề nghiệp cha",
    "Họ tên mẹ",
    "Năm sinh mẹ",
    "Số điện thoại mẹ",
    "Nghề nghiệp mẹ",
    "Họ tên quan hệ khác",
    "Năm sinh quan hệ khác",
    "Số điện thoại quan hệ khác",
    "Nghề nghiệp quan hệ khác"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# The row used to run out to AC1 (29 cols); now it only runs to AA1 (27 cols) -
# clear the two trailing cells so the sheet's used range / dimension shrinks back.
$ws.Range("AB1:AC1").Clear() | Out-Null

# Column widths follow the header text (bestFit) in the source file; re-apply the
# closest widths for the new column order. Column A ("STT") keeps its original
# width untouched since it didn't change. (Values below are the ColumnWidth
# inputs that reproduce the target stored widths through this host's rounding.)
$widths = @(
    7.166666666666667,
    3.5,
    3.0,
    2.3333333333333335,
    3.0,
    7.666666666666667,
    7.0,
    6.333333333333333,
    13.833333333333334,
    10.666666666666666,
    5.666666666666667,
    3.1666666666666665,
    2.6666666666666665,
    11.166666666666666,
    8.833333333333334,
    11.0,
    14.166666666666666,
    13.833333333333334,
    8.333333333333334,
    10.666666666666666,
    13.833333333333334,
    13.333333333333334,
    17.0,
    19.166666666666668,
    22.333333333333332,
    22.0
)
for ($i = 0; $i -lt $widths.Length; $i++) {
    $ws.Columns.Item($i + 2).ColumnWidth = $widths[$i]
}

# Selection moved from H11 to K10 in the saved view.
$ws.Range("K10").Select() | Out-Null

Write-Output "done"
